# Update the "Presentations" slide (slide 2) list of AWS slide-deck links:
#   - "Amazon EC2 – check back here in 5-7 days" becomes
#     "Amazon EC2 – Coming Soon"
#   - a new paragraph/hyperlink for "Amazon Relational Database Service (RDS)"
#     is appended right after it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("Content Placeholder 2")
$tr = $shp.TextFrame.TextRange

$full = $tr.Text
$needle = "Amazon EC2"
$idx = $full.IndexOf($needle)
$ec2Len = $full.Length - $idx
$ec2Range = $tr.Characters($idx + 1, $ec2Len)
$ec2Range.Text = "Amazon EC2 – Coming Soon"

$newText = "Amazon Relational Database Service (RDS)"
$beforeLen = $tr.Length
$null = $tr.InsertAfter("`r" + $newText)

$rdsRange = $tr.Characters($beforeLen + 2, $newText.Length)
$rdsRange.ActionSettings.Item(1).Hyperlink.Address = "https://s3.us-east-2.amazonaws.com/public.jeff-anderson.com/IntroToAWS-RDS.pdf"
